$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
Write-Host ($excel.ActiveWindow | Get-Member -MemberType Property | Out-String)
